$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update recalculated statistic values in existing rows (rows 2-133) ---
$ws.Range("G2").Value = 1.95852704099881
$ws.Range("I2").Value = 3.60088
$ws.Range("N2").Value = 2.99354
$ws.Range("G9").Value = 0.0176341983931321
$ws.Range("G10").Value = 0.0176341983931321
$ws.Range("G11").Value = 0.719032147324523
$ws.Range("G12").Value = 0.719032147324523
$ws.Range("G13").Value = 0.7414666666666671
$ws.Range("L13").Value = 0.1468
$ws.Range("M13").Value = 1.35397
$ws.Range("N13").Value = 1.56126
$ws.Range("G14").Value = 0.7414666666666671
$ws.Range("L14").Value = 0.1468
$ws.Range("M14").Value = 1.35397
$ws.Range("N14").Value = 1.56126
$ws.Range("G26").Value = 0.0152099621428922
$ws.Range("G27").Value = 0.0152099621428922
$ws.Range("G28").Value = 0.691211245045376
$ws.Range("G29").Value = 0.691211245045376
$ws.Range("G30").Value = 0.711534652507853
$ws.Range("L30").Value = 0.1444
$ws.Range("M30").Value = 1.31939
$ws.Range("N30").Value = 1.54868
$ws.Range("G31").Value = 0.711534652507853
$ws.Range("L31").Value = 0.1444
$ws.Range("M31").Value = 1.31939
$ws.Range("N31").Value = 1.54868
$ws.Range("F43").Value = 0.00277
$ws.Range("G43").Value = 0.007949051368473599
$ws.Range("L43").Value = 0.00117
$ws.Range("F44").Value = 0.00277
$ws.Range("G44").Value = 0.007949051368473599
$ws.Range("L44").Value = 0.00117
$ws.Range("G45").Value = 0.73269457837871
$ws.Range("G46").Value = 0.73269457837871
$ws.Range("F47").Value = 0.53415
$ws.Range("G47").Value = 0.74834131917452
$ws.Range("M47").Value = 1.36882
$ws.Range("F48").Value = 0.53415
$ws.Range("G48").Value = 0.74834131917452
$ws.Range("M48").Value = 1.36882
$ws.Range("G56").Value = 584.589648228522
$ws.Range("H56").Value = 5199
$ws.Range("G57").Value = 584.589648228522
$ws.Range("H57").Value = 5199
$ws.Range("G58").Value = 584.589648228522
$ws.Range("H58").Value = 5199
$ws.Range("G59").Value = 584.589648228522
$ws.Range("H59").Value = 5199
$ws.Range("F60").Value = 0.00341
$ws.Range("G60").Value = 0.0069353679948533
$ws.Range("L60").Value = 0.00227
$ws.Range("F61").Value = 0.00341
$ws.Range("G61").Value = 0.0069353679948533
$ws.Range("L61").Value = 0.00227
$ws.Range("G62").Value = 0.761407911712043
$ws.Range("G63").Value = 0.761407911712043
$ws.Range("G64").Value = 0.7780296525078531
$ws.Range("M64").Value = 1.36882
$ws.Range("G65").Value = 0.7780296525078531
$ws.Range("M65").Value = 1.36882
$ws.Range("G73").Value = 500.572981561855
$ws.Range("H73").Value = 4573.37889371131
$ws.Range("G74").Value = 500.572981561855
$ws.Range("H74").Value = 4573.37889371131
$ws.Range("G75").Value = 500.572981561855
$ws.Range("H75").Value = 4573.37889371131
$ws.Range("G76").Value = 500.572981561855
$ws.Range("H76").Value = 4573.37889371131
$ws.Range("F77").Value = 0.00484
$ws.Range("G77").Value = 0.009184332373738001
$ws.Range("I77").Value = 0.03008
$ws.Range("L77").Value = 0.00398
$ws.Range("F78").Value = 0.00484
$ws.Range("G78").Value = 0.009184332373738001
$ws.Range("I78").Value = 0.03008
$ws.Range("L78").Value = 0.00398
$ws.Range("G79").Value = 0.783702911712043
$ws.Range("G80").Value = 0.783702911712043
$ws.Range("G81").Value = 0.80253131917452
$ws.Range("G82").Value = 0.80253131917452
$ws.Range("G90").Value = 439.514896503582
$ws.Range("H90").Value = 4573.37889371131
$ws.Range("G91").Value = 439.514896503582
$ws.Range("H91").Value = 4573.37889371131
$ws.Range("G92").Value = 439.514896503582
$ws.Range("H92").Value = 4573.37889371131
$ws.Range("G93").Value = 439.514896503582
$ws.Range("H93").Value = 4573.37889371131
$ws.Range("G94").Value = 0.0093047380689288
$ws.Range("I94").Value = 0.03009
$ws.Range("L94").Value = 0.00412
$ws.Range("N94").Value = 0.02288
$ws.Range("G95").Value = 0.0093047380689288
$ws.Range("I95").Value = 0.03009
$ws.Range("L95").Value = 0.00412
$ws.Range("N95").Value = 0.02288
$ws.Range("G96").Value = 0.761362350012067
$ws.Range("G97").Value = 0.761362350012067
$ws.Range("G98").Value = 0.779662306606473
$ws.Range("G99").Value = 0.779662306606473
$ws.Range("G107").Value = 414.887777859514
$ws.Range("H107").Value = 4573.37889371131
$ws.Range("G108").Value = 414.887777859514
$ws.Range("H108").Value = 4573.37889371131
$ws.Range("G109").Value = 414.887777859514
$ws.Range("H109").Value = 4573.37889371131
$ws.Range("G110").Value = 414.887777859514
$ws.Range("H110").Value = 4573.37889371131
$ws.Range("F111").Value = 0.00468
$ws.Range("G111").Value = 0.0097790026745066
$ws.Range("I111").Value = 0.04701
$ws.Range("L111").Value = 0.00456
$ws.Range("N111").Value = 0.02288
$ws.Range("F112").Value = 0.00468
$ws.Range("G112").Value = 0.0097790026745066
$ws.Range("I112").Value = 0.04701
$ws.Range("L112").Value = 0.00456
$ws.Range("N112").Value = 0.02288
$ws.Range("G113").Value = 0.724839599076054
$ws.Range("G114").Value = 0.724839599076054
$ws.Range("G115").Value = 0.743145710835775
$ws.Range("G116").Value = 0.743145710835775
$ws.Range("G124").Value = 482.427699889672
$ws.Range("H124").Value = 4573.37889371131
$ws.Range("G125").Value = 482.427699889672
$ws.Range("H125").Value = 4573.37889371131
$ws.Range("G126").Value = 482.427699889672
$ws.Range("H126").Value = 4573.37889371131
$ws.Range("G127").Value = 482.427699889672
$ws.Range("H127").Value = 4573.37889371131
$ws.Range("G128").Value = 0.0103261462869458
$ws.Range("N128").Value = 0.02798
$ws.Range("G129").Value = 0.0103261462869458
$ws.Range("N129").Value = 0.02798
$ws.Range("G130").Value = 0.588869058692758
$ws.Range("G131").Value = 0.588869058692758
$ws.Range("G132").Value = 0.6073087182335209
$ws.Range("G133").Value = 0.6073087182335209

# --- 2. Append new rows 138-154 (2019 - 2023 results) ---
# Row 138: Visual Clarity (Sediment class 4)
$ws.Range("A138").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B138").Value = "Visual Clarity (Sediment class 4)"
$ws.Range("C138").Value = "A"
$ws.Range("D138").Value = "2019 - 2023"
$ws.Range("E138").Value = "RepSite"
$ws.Range("F138").Value = 2.3
$ws.Range("G138").Value = 2.31538461538462
$ws.Range("H138").Value = 4.6
$ws.Range("I138").Value = 4.335
$ws.Range("J138").Value = "'"
$ws.Range("J138").Style = "Normal"
$ws.Range("K138").Value = "'"
$ws.Range("K138").Style = "Normal"
$ws.Range("L138").Value = 2.75
$ws.Range("M138").Value = 3.283
$ws.Range("N138").Value = 4.068
$ws.Range("O138").Value = 1878758.348
$ws.Range("P138").Value = 5554292.094
$ws.Range("Q138").Value = "Tararua District"
$ws.Range("R138").Value = "Manawatū"
$ws.Range("S138").Value = "Upper Manawatu"
$ws.Range("T138").Value = "Mana_1a"
$ws.Range("U138").Value = "m"

# Row 139: DRP (95th Percentile)
$ws.Range("A139").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B139").Value = "DRP (95th Percentile)"
$ws.Range("C139").Value = "D"
$ws.Range("D139").Value = "2019 - 2023"
$ws.Range("E139").Value = "RepSite"
$ws.Range("F139").Value = 0.032
$ws.Range("G139").Value = 0.034859649122807
$ws.Range("H139").Value = 0.069
$ws.Range("I139").Value = 0.0663
$ws.Range("J139").Value = "'"
$ws.Range("J139").Style = "Normal"
$ws.Range("K139").Value = "'"
$ws.Range("K139").Style = "Normal"
$ws.Range("L139").Value = 0.039
$ws.Range("M139").Value = 0.05381
$ws.Range("N139").Value = 0.06082
$ws.Range("O139").Value = 1878758.348
$ws.Range("P139").Value = 5554292.094
$ws.Range("Q139").Value = "Tararua District"
$ws.Range("R139").Value = "Manawatū"
$ws.Range("S139").Value = "Upper Manawatu"
$ws.Range("T139").Value = "Mana_1a"
$ws.Range("U139").Value = "mg/L"

# Row 140: DRP (Median)
$ws.Range("A140").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B140").Value = "DRP (Median)"
$ws.Range("C140").Value = "D"
$ws.Range("D140").Value = "2019 - 2023"
$ws.Range("E140").Value = "RepSite"
$ws.Range("F140").Value = 0.032
$ws.Range("G140").Value = 0.034859649122807
$ws.Range("H140").Value = 0.069
$ws.Range("I140").Value = 0.0663
$ws.Range("J140").Value = "'"
$ws.Range("J140").Style = "Normal"
$ws.Range("K140").Value = "'"
$ws.Range("K140").Style = "Normal"
$ws.Range("L140").Value = 0.039
$ws.Range("M140").Value = 0.05381
$ws.Range("N140").Value = 0.06082
$ws.Range("O140").Value = 1878758.348
$ws.Range("P140").Value = 5554292.094
$ws.Range("Q140").Value = "Tararua District"
$ws.Range("R140").Value = "Manawatū"
$ws.Range("S140").Value = "Upper Manawatu"
$ws.Range("T140").Value = "Mana_1a"
$ws.Range("U140").Value = "mg/L"

# Row 141: E coli (>260)
$ws.Range("A141").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B141").Value = "E coli (>260)"
$ws.Range("C141").Value = "D"
$ws.Range("D141").Value = "2019 - 2023"
$ws.Range("E141").Value = "RepSite"
$ws.Range("F141").Value = 187
$ws.Range("G141").Value = 443.298245614035
$ws.Range("H141").Value = 3300
$ws.Range("I141").Value = 2265
$ws.Range("J141").Value = 17.5438596491228
$ws.Range("K141").Value = 36.8421052631579
$ws.Range("L141").Value = 210
$ws.Range("M141").Value = 612.56
$ws.Range("N141").Value = 1670
$ws.Range("O141").Value = 1878758.348
$ws.Range("P141").Value = 5554292.094
$ws.Range("Q141").Value = "Tararua District"
$ws.Range("R141").Value = "Manawatū"
$ws.Range("S141").Value = "Upper Manawatu"
$ws.Range("T141").Value = "Mana_1a"
$ws.Range("U141").Value = "% exceedances over 260/100 mL"

# Row 142: E coli (>540)
$ws.Range("A142").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B142").Value = "E coli (>540)"
$ws.Range("C142").Value = "C"
$ws.Range("D142").Value = "2019 - 2023"
$ws.Range("E142").Value = "RepSite"
$ws.Range("F142").Value = 187
$ws.Range("G142").Value = 443.298245614035
$ws.Range("H142").Value = 3300
$ws.Range("I142").Value = 2265
$ws.Range("J142").Value = 17.5438596491228
$ws.Range("K142").Value = 36.8421052631579
$ws.Range("L142").Value = 210
$ws.Range("M142").Value = 612.56
$ws.Range("N142").Value = 1670
$ws.Range("O142").Value = 1878758.348
$ws.Range("P142").Value = 5554292.094
$ws.Range("Q142").Value = "Tararua District"
$ws.Range("R142").Value = "Manawatū"
$ws.Range("S142").Value = "Upper Manawatu"
$ws.Range("T142").Value = "Mana_1a"
$ws.Range("U142").Value = "% exceedances over 540/100 mL"

# Row 143: E coli (Median)
$ws.Range("A143").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B143").Value = "E coli (Median)"
$ws.Range("C143").Value = "D"
$ws.Range("D143").Value = "2019 - 2023"
$ws.Range("E143").Value = "RepSite"
$ws.Range("F143").Value = 187
$ws.Range("G143").Value = 443.298245614035
$ws.Range("H143").Value = 3300
$ws.Range("I143").Value = 2265
$ws.Range("J143").Value = 17.5438596491228
$ws.Range("K143").Value = 36.8421052631579
$ws.Range("L143").Value = 210
$ws.Range("M143").Value = 612.56
$ws.Range("N143").Value = 1670
$ws.Range("O143").Value = 1878758.348
$ws.Range("P143").Value = 5554292.094
$ws.Range("Q143").Value = "Tararua District"
$ws.Range("R143").Value = "Manawatū"
$ws.Range("S143").Value = "Upper Manawatu"
$ws.Range("T143").Value = "Mana_1a"
$ws.Range("U143").Value = "E. coli/100 mL"

# Row 144: E coli (95th Percentile)
$ws.Range("A144").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B144").Value = "E coli (95th Percentile)"
$ws.Range("C144").Value = "E"
$ws.Range("D144").Value = "2019 - 2023"
$ws.Range("E144").Value = "RepSite"
$ws.Range("F144").Value = 187
$ws.Range("G144").Value = 443.298245614035
$ws.Range("H144").Value = 3300
$ws.Range("I144").Value = 2265
$ws.Range("J144").Value = 17.5438596491228
$ws.Range("K144").Value = 36.8421052631579
$ws.Range("L144").Value = 210
$ws.Range("M144").Value = 612.56
$ws.Range("N144").Value = 1670
$ws.Range("O144").Value = 1878758.348
$ws.Range("P144").Value = 5554292.094
$ws.Range("Q144").Value = "Tararua District"
$ws.Range("R144").Value = "Manawatū"
$ws.Range("S144").Value = "Upper Manawatu"
$ws.Range("T144").Value = "Mana_1a"
$ws.Range("U144").Value = "E. coli/100 mL"

# Row 145: Ammoniacal-N (95th Percentile)
$ws.Range("A145").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B145").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C145").Value = "A"
$ws.Range("D145").Value = "2019 - 2023"
$ws.Range("E145").Value = "RepSite"
$ws.Range("F145").Value = 0.00685
$ws.Range("G145").Value = 0.0102605221859975
$ws.Range("H145").Value = 0.0730536907854253
$ws.Range("I145").Value = 0.02965
$ws.Range("J145").Value = "'"
$ws.Range("J145").Style = "Normal"
$ws.Range("K145").Value = "'"
$ws.Range("K145").Style = "Normal"
$ws.Range("L145").Value = 0.00559
$ws.Range("M145").Value = 0.01467
$ws.Range("N145").Value = 0.02585
$ws.Range("O145").Value = 1878758.348
$ws.Range("P145").Value = 5554292.094
$ws.Range("Q145").Value = "Tararua District"
$ws.Range("R145").Value = "Manawatū"
$ws.Range("S145").Value = "Upper Manawatu"
$ws.Range("T145").Value = "Mana_1a"
$ws.Range("U145").Value = "mg NH4-N/L"

# Row 146: Ammoniacal-N (Median)
$ws.Range("A146").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B146").Value = "Ammoniacal-N (Median)"
$ws.Range("C146").Value = "A"
$ws.Range("D146").Value = "2019 - 2023"
$ws.Range("E146").Value = "RepSite"
$ws.Range("F146").Value = 0.00685
$ws.Range("G146").Value = 0.0102605221859975
$ws.Range("H146").Value = 0.0730536907854253
$ws.Range("I146").Value = 0.02965
$ws.Range("J146").Value = "'"
$ws.Range("J146").Style = "Normal"
$ws.Range("K146").Value = "'"
$ws.Range("K146").Style = "Normal"
$ws.Range("L146").Value = 0.00559
$ws.Range("M146").Value = 0.01467
$ws.Range("N146").Value = 0.02585
$ws.Range("O146").Value = 1878758.348
$ws.Range("P146").Value = 5554292.094
$ws.Range("Q146").Value = "Tararua District"
$ws.Range("R146").Value = "Manawatū"
$ws.Range("S146").Value = "Upper Manawatu"
$ws.Range("T146").Value = "Mana_1a"
$ws.Range("U146").Value = "mg NH4-N/L"

# Row 147: Nitrate-N (95th Percentile)
$ws.Range("A147").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B147").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C147").Value = "B"
$ws.Range("D147").Value = "2019 - 2023"
$ws.Range("E147").Value = "RepSite"
$ws.Range("F147").Value = 0.495
$ws.Range("G147").Value = 0.612921690271705
$ws.Range("H147").Value = 2.14
$ws.Range("I147").Value = 1.743
$ws.Range("J147").Value = "'"
$ws.Range("J147").Style = "Normal"
$ws.Range("K147").Value = "'"
$ws.Range("K147").Style = "Normal"
$ws.Range("L147").Value = 0.126
$ws.Range("M147").Value = 0.9933
$ws.Range("N147").Value = 1.5422
$ws.Range("O147").Value = 1878758.348
$ws.Range("P147").Value = 5554292.094
$ws.Range("Q147").Value = "Tararua District"
$ws.Range("R147").Value = "Manawatū"
$ws.Range("S147").Value = "Upper Manawatu"
$ws.Range("T147").Value = "Mana_1a"
$ws.Range("U147").Value = "mg NO3-N/L"

# Row 148: Nitrate-N (Median)
$ws.Range("A148").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B148").Value = "Nitrate-N (Median)"
$ws.Range("C148").Value = "A"
$ws.Range("D148").Value = "2019 - 2023"
$ws.Range("E148").Value = "RepSite"
$ws.Range("F148").Value = 0.495
$ws.Range("G148").Value = 0.612921690271705
$ws.Range("H148").Value = 2.14
$ws.Range("I148").Value = 1.743
$ws.Range("J148").Value = "'"
$ws.Range("J148").Style = "Normal"
$ws.Range("K148").Value = "'"
$ws.Range("K148").Style = "Normal"
$ws.Range("L148").Value = 0.126
$ws.Range("M148").Value = 0.9933
$ws.Range("N148").Value = 1.5422
$ws.Range("O148").Value = 1878758.348
$ws.Range("P148").Value = 5554292.094
$ws.Range("Q148").Value = "Tararua District"
$ws.Range("R148").Value = "Manawatū"
$ws.Range("S148").Value = "Upper Manawatu"
$ws.Range("T148").Value = "Mana_1a"
$ws.Range("U148").Value = "mg NO3-N/L"

# Row 149: Soluble Inorganic Nitrogen (95th Percentile)
$ws.Range("A149").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B149").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("C149").Value = "'"
$ws.Range("C149").Style = "Normal"
$ws.Range("D149").Value = "2019 - 2023"
$ws.Range("E149").Value = "RepSite"
$ws.Range("F149").Value = 0.515
$ws.Range("G149").Value = 0.630817490163346
$ws.Range("H149").Value = 2.155
$ws.Range("I149").Value = 1.7695
$ws.Range("J149").Value = "'"
$ws.Range("J149").Style = "Normal"
$ws.Range("K149").Value = "'"
$ws.Range("K149").Style = "Normal"
$ws.Range("L149").Value = 0.142
$ws.Range("M149").Value = 1.01316
$ws.Range("N149").Value = 1.56878
$ws.Range("O149").Value = 1878758.348
$ws.Range("P149").Value = 5554292.094
$ws.Range("Q149").Value = "Tararua District"
$ws.Range("R149").Value = "Manawatū"
$ws.Range("S149").Value = "Upper Manawatu"
$ws.Range("T149").Value = "Mana_1a"
$ws.Range("U149").Value = "g/m3"

# Row 150: Soluble Inorganic Nitrogen (Median)
$ws.Range("A150").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B150").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("C150").Value = "'"
$ws.Range("C150").Style = "Normal"
$ws.Range("D150").Value = "2019 - 2023"
$ws.Range("E150").Value = "RepSite"
$ws.Range("F150").Value = 0.515
$ws.Range("G150").Value = 0.630817490163346
$ws.Range("H150").Value = 2.155
$ws.Range("I150").Value = 1.7695
$ws.Range("J150").Value = "'"
$ws.Range("J150").Style = "Normal"
$ws.Range("K150").Value = "'"
$ws.Range("K150").Style = "Normal"
$ws.Range("L150").Value = 0.142
$ws.Range("M150").Value = 1.01316
$ws.Range("N150").Value = 1.56878
$ws.Range("O150").Value = 1878758.348
$ws.Range("P150").Value = 5554292.094
$ws.Range("Q150").Value = "Tararua District"
$ws.Range("R150").Value = "Manawatū"
$ws.Range("S150").Value = "Upper Manawatu"
$ws.Range("T150").Value = "Mana_1a"
$ws.Range("U150").Value = "g/m3"

# Row 151: Total Nitrogen (95th Percentile)
$ws.Range("A151").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B151").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("C151").Value = "'"
$ws.Range("C151").Style = "Normal"
$ws.Range("D151").Value = "2019 - 2023"
$ws.Range("E151").Value = "RepSite"
$ws.Range("F151").Value = 0.78
$ws.Range("G151").Value = 0.862456140350877
$ws.Range("H151").Value = 2.53
$ws.Range("I151").Value = 1.9655
$ws.Range("J151").Value = "'"
$ws.Range("J151").Style = "Normal"
$ws.Range("K151").Value = "'"
$ws.Range("K151").Style = "Normal"
$ws.Range("L151").Value = 0.34
$ws.Range("M151").Value = 1.3429
$ws.Range("N151").Value = 1.7476
$ws.Range("O151").Value = 1878758.348
$ws.Range("P151").Value = 5554292.094
$ws.Range("Q151").Value = "Tararua District"
$ws.Range("R151").Value = "Manawatū"
$ws.Range("S151").Value = "Upper Manawatu"
$ws.Range("T151").Value = "Mana_1a"
$ws.Range("U151").Value = "g/m3"

# Row 152: Total Nitrogen (Median)
$ws.Range("A152").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B152").Value = "Total Nitrogen (Median)"
$ws.Range("C152").Value = "'"
$ws.Range("C152").Style = "Normal"
$ws.Range("D152").Value = "2019 - 2023"
$ws.Range("E152").Value = "RepSite"
$ws.Range("F152").Value = 0.78
$ws.Range("G152").Value = 0.862456140350877
$ws.Range("H152").Value = 2.53
$ws.Range("I152").Value = 1.9655
$ws.Range("J152").Value = "'"
$ws.Range("J152").Style = "Normal"
$ws.Range("K152").Value = "'"
$ws.Range("K152").Style = "Normal"
$ws.Range("L152").Value = 0.34
$ws.Range("M152").Value = 1.3429
$ws.Range("N152").Value = 1.7476
$ws.Range("O152").Value = 1878758.348
$ws.Range("P152").Value = 5554292.094
$ws.Range("Q152").Value = "Tararua District"
$ws.Range("R152").Value = "Manawatū"
$ws.Range("S152").Value = "Upper Manawatu"
$ws.Range("T152").Value = "Mana_1a"
$ws.Range("U152").Value = "g/m3"

# Row 153: Total Phosphorus (95th Percentile)
$ws.Range("A153").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B153").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("C153").Value = "'"
$ws.Range("C153").Style = "Normal"
$ws.Range("D153").Value = "2019 - 2023"
$ws.Range("E153").Value = "RepSite"
$ws.Range("F153").Value = 0.049
$ws.Range("G153").Value = 0.0511754385964912
$ws.Range("H153").Value = 0.181
$ws.Range("I153").Value = 0.0769
$ws.Range("J153").Value = "'"
$ws.Range("J153").Style = "Normal"
$ws.Range("K153").Value = "'"
$ws.Range("K153").Style = "Normal"
$ws.Range("L153").Value = 0.055
$ws.Range("M153").Value = 0.065
$ws.Range("N153").Value = 0.06994
$ws.Range("O153").Value = 1878758.348
$ws.Range("P153").Value = 5554292.094
$ws.Range("Q153").Value = "Tararua District"
$ws.Range("R153").Value = "Manawatū"
$ws.Range("S153").Value = "Upper Manawatu"
$ws.Range("T153").Value = "Mana_1a"
$ws.Range("U153").Value = "g/m3"

# Row 154: Total Phosphorus (Median)
$ws.Range("A154").Value = "Mangarangiora at u/s Ormondville STP"
$ws.Range("B154").Value = "Total Phosphorus (Median)"
$ws.Range("C154").Value = "'"
$ws.Range("C154").Style = "Normal"
$ws.Range("D154").Value = "2019 - 2023"
$ws.Range("E154").Value = "RepSite"
$ws.Range("F154").Value = 0.049
$ws.Range("G154").Value = 0.0511754385964912
$ws.Range("H154").Value = 0.181
$ws.Range("I154").Value = 0.0769
$ws.Range("J154").Value = "'"
$ws.Range("J154").Style = "Normal"
$ws.Range("K154").Value = "'"
$ws.Range("K154").Style = "Normal"
$ws.Range("L154").Value = 0.055
$ws.Range("M154").Value = 0.065
$ws.Range("N154").Value = 0.06994
$ws.Range("O154").Value = 1878758.348
$ws.Range("P154").Value = 5554292.094
$ws.Range("Q154").Value = "Tararua District"
$ws.Range("R154").Value = "Manawatū"
$ws.Range("S154").Value = "Upper Manawatu"
$ws.Range("T154").Value = "Mana_1a"
$ws.Range("U154").Value = "g/m3"

